$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newRow = 33

$ws.Cells.Item($newRow, 1).Value = "The Use of Standard Parenting Management Training in Addressing Disruptive Mood Dysregulation Disorder: A Pilot Study"
$ws.Cells.Item($newRow, 2).Value = "Gary Byrne, Graham Connon"
$ws.Cells.Item($newRow, 3).Value = "Journal of Contemporary Psychotherapy"
$ws.Cells.Item($newRow, 4).Value = 2021
$ws.Cells.Item($newRow, 5).Value = "https://doi.org/10.1007/s10879-021-09489-5"
$ws.Cells.Item($newRow, 6).Value = "Parent management training has demonstrated effectiveness in the treatment of child behavioural issues and associated conduct difficulties. Anger, aggression, and irritability are common symptoms amongst children presenting with disruptive mood dysregulation disorder. Currently, there are no well-established evidence-based interventions for children presenting with symptoms of disruptive mood dysregulation disorder. This pilot study aims to assess if a standard, well-established, parent management training program (group Triple P) may be effective in addressing disruptive mood dysregulation disorder symptoms. Thirteen parents of children who presented with disruptive mood dysregulation disorder or subthreshold symptoms completed the Triple P behavioural management program (Level 4). Post-treatment, parents reported no significant change on childhood irritability. However, parents noted significant improvement on child overt aggression, behavioural difficulties and an increase in child pro-social behaviours. Despite the many limitations inherent in this pilot study, results suggest that standard parent management training may be useful in addressing overt aggression but not irritability."
$ws.Cells.Item($newRow, 7).Value = "ARI used as part of measures"
$ws.Cells.Item($newRow, 8).Value = "Treatment"

# Column I (Type_Secondary) is left blank for this row, matching the
# existing blank placeholder cells used throughout the table (e.g. I3, I4).
# Touching the style (without setting a Value) materialises the cell
# reference in the sheet without giving it any content or formatting.
$ws.Cells.Item($newRow, 9).Style = "Normal"
